$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.341.78"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "1.873.09"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.61"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4682"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06418"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07789"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.888.59"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.91"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7189"
$ws.Range("E14").Value = "  -3.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.130"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.38"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "30.328.66"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.139.77"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007377"
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.209"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.236"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.68"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.004"
$ws.Range("E26").Value = "  -2.04%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  -1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.331"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09574"
$ws.Range("E30").Value = "  -3.75%  "
$ws.Range("E31").Value = "  -2.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.202"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.067"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04803"
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.115"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6841"
$ws.Range("E36").Value = "  -1.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01870"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.807"
$ws.Range("E39").Value = "  +1.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.219"
$ws.Range("E40").Value = "  -0.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.12"
$ws.Range("E41").Value = "  +1.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.929"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4211"
$ws.Range("E43").Value = "  +1.24%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8220"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.50"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.580"
$ws.Range("E47").Value = "  +2.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.97"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.891"
$ws.Range("E49").Value = "  -1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "895.80"
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05726"
$ws.Range("E51").Value = "  +1.08%  "
